# "alignment tree and aligners. various corrections."
#
# The feature-location table for BTV segments 9 and 10 had two stray
# duplicate rows (leftover from a bad paste of the 5'UTR boundaries) and
# one mislabeled product name. Clean those up:
#   - Segment 9 block: drop the extra "Seg-9_5UTR 1 27" row (row 35) that
#     duplicated the segment-9 5'UTR entry already on row 34.
#   - Segment 8 product name I31 was mislabeled "VS2"; it should read "NS2".
#   - Segment 10 block: drop the extra "Seg-10_5UTR 1 58" row that
#     duplicated the segment-10 5'UTR entry already on row 40 (=old row 41).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the mislabeled product name on segment 8's CDS row.
$ws.Range("I31").Value = "NS2"

# Remove the stray duplicate "Seg-9_5UTR" row under segment 9.
$ws.Rows("35").Delete()

# Remove the stray duplicate "Seg-10_5UTR" row under segment 10 (this row
# shifted up to 41 after the first deletion above).
$ws.Rows("41").Delete()

# Restore the view: zoomed in on the segment-10 block the edit focused on.
$ws.Select()
$excel.ActiveWindow.Zoom = 200
$ws.Range("G22:I24").Select()
